$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to text format so numeric-looking strings
# (e.g. "1.00", "91.872.05") are preserved exactly as text, matching the source data.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '91.872.05'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '3.129.56'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '242.33'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = '618.45'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('D7').Value = '1.10'
$ws.Range('E7').Value = '  -5.50%  '
$ws.Range('D8').Value = '0.388'
$ws.Range('E8').Value = '  +3.98%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '3.126.27'
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = '0.0000252'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '35.34'
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = '5.61'
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').Value = '91.527.52'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').Value = '3.711.10'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '3.116.57'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('D20').Value = '15.00'
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('D21').Value = '5.92'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').Value = '457.56'
$ws.Range('E22').Value = '  +1.27%  '
$ws.Range('D23').Value = '0.0000202'
$ws.Range('E23').Value = '  -5.56%  '
$ws.Range('D24').Value = '9.25'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('D25').Value = '5.96'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '89.33'
$ws.Range('E26').Value = '  -4.34%  '
$ws.Range('D27').Value = '11.74'
$ws.Range('E27').Value = '  -2.11%  '
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('E29').Value = '  +15.37%  '
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').Value = '0.168'
$ws.Range('E31').Value = '  -6.10%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.226'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').Value = '9.45'
$ws.Range('E33').Value = '  +3.55%  '
$ws.Range('D34').Value = '0.177'
$ws.Range('E34').Value = '  +7.82%  '
$ws.Range('D35').Value = '26.52'
$ws.Range('E35').Value = '  -1.81%  '
$ws.Range('E36').Value = '  -2.19%  '
$ws.Range('E37').Value = '  +1.08%  '
$ws.Range('D38').Value = '491.80'
$ws.Range('E38').Value = '  -2.29%  '
$ws.Range('D39').Value = '3.88'
$ws.Range('E39').Value = '  -7.25%  '
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('D41').Value = '0.440'
$ws.Range('E41').Value = '  +3.22%  '
$ws.Range('D42').Value = '3.39'
$ws.Range('E42').Value = '  -6.49%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E45').Value = '  -28.96%  '
$ws.Range('D46').Value = '0.711'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('D47').Value = '156.61'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').Value = '1.92'
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('D49').Value = '1.35'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('D50').Value = '4.48'
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('D51').Value = '0.0327'
$ws.Range('E51').Value = '  +0.14%  '

# Restore the default (Normal) style on the Price column now that values are set,
# so no residual text-format style is left applied to the cells.
$priceRange.Style = "Normal"
